$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 84
$ws.Cells.Item(84,1).Value = "2025-10-20 18:33:13"
$ws.Cells.Item(84,2).Value = "Noahs channel"
$ws.Cells.Item(84,4).Value = "Unknown"
$ws.Cells.Item(84,5).Value = "Test message"
$ws.Cells.Item(84,7).Value = "Noahs channel"

# Row 85
$ws.Cells.Item(85,1).Value = "2025-10-20 18:34:36"
$ws.Cells.Item(85,2).Value = "Noahs channel"
$ws.Cells.Item(85,4).Value = "Unknown"
$ws.Cells.Item(85,5).Value = "Test message"
$ws.Cells.Item(85,7).Value = "Noahs channel"

# Row 86
$ws.Cells.Item(86,1).Value = "2025-10-20 18:35:28"
$ws.Cells.Item(86,2).Value = "Noah Dubitzky"
$ws.Cells.Item(86,3).Value = 8450689526
$ws.Cells.Item(86,4).NumberFormat = "@"
$ws.Cells.Item(86,4).Value = "13052054965"
$ws.Cells.Item(86,4).Style = "Normal"
$ws.Cells.Item(86,5).Value = "Hey man"

# Row 87
$ws.Cells.Item(87,1).Value = "2025-10-20 18:35:43"
$ws.Cells.Item(87,2).Value = "Noahs channel"
$ws.Cells.Item(87,4).Value = "Unknown"
$ws.Cells.Item(87,5).Value = "Test message"
$ws.Cells.Item(87,7).Value = "Noahs channel"

# Row 88
$ws.Cells.Item(88,1).Value = "2025-10-20 18:42:13"
$ws.Cells.Item(88,2).Value = "Noahs channel"
$ws.Cells.Item(88,4).Value = "Unknown"
$ws.Cells.Item(88,5).Value = "Hey guys"
$ws.Cells.Item(88,7).Value = "Noahs channel"

# Row 89
$ws.Cells.Item(89,1).Value = "2025-10-20 18:45:23"
$ws.Cells.Item(89,2).Value = "Noahs channel"
$ws.Cells.Item(89,4).Value = "Unknown"
$ws.Cells.Item(89,5).Value = "Test"
$ws.Cells.Item(89,7).Value = "Noahs channel"

# Row 90
$ws.Cells.Item(90,1).Value = "2025-10-20 19:00:08"
$ws.Cells.Item(90,2).Value = "Noahs channel"
$ws.Cells.Item(90,4).Value = "Unknown"
$ws.Cells.Item(90,5).Value = "Test"
$ws.Cells.Item(90,7).Value = "Noahs channel"

# Row 91
$ws.Cells.Item(91,1).Value = "2025-10-20 19:01:49"
$ws.Cells.Item(91,2).Value = "Noahs channel"
$ws.Cells.Item(91,4).Value = "Unknown"
$ws.Cells.Item(91,5).Value = "Test"
$ws.Cells.Item(91,7).Value = "Noahs channel"

# Row 92
$ws.Cells.Item(92,1).Value = "2025-10-20 19:02:54"
$ws.Cells.Item(92,2).Value = "Noahs channel"
$ws.Cells.Item(92,4).Value = "Unknown"
$ws.Cells.Item(92,5).Value = "Test"
$ws.Cells.Item(92,7).Value = "Noahs channel"

# Row 93
$ws.Cells.Item(93,1).Value = "2025-10-20 19:08:35"
$ws.Cells.Item(93,2).Value = "Noahs channel"
$ws.Cells.Item(93,4).Value = "Unknown"
$ws.Cells.Item(93,5).Value = "Test message"
$ws.Cells.Item(93,7).Value = "Noahs channel"

# Row 94
$ws.Cells.Item(94,1).Value = "2025-10-20 19:08:35"
$ws.Cells.Item(94,2).Value = "Noahs channel"
$ws.Cells.Item(94,4).Value = "Unknown"
$ws.Cells.Item(94,5).Value = "Test message"
$ws.Cells.Item(94,7).Value = "Noahs channel"

# Row 95
$ws.Cells.Item(95,1).Value = "2025-10-20 19:09:30"
$ws.Cells.Item(95,2).Value = "Noahs channel"
$ws.Cells.Item(95,4).Value = "Unknown"
$ws.Cells.Item(95,5).Value = "Test"
$ws.Cells.Item(95,7).Value = "Noahs channel"

# Row 96
$ws.Cells.Item(96,1).Value = "2025-10-20 19:10:11"
$ws.Cells.Item(96,2).Value = "Noah Dubitzky"
$ws.Cells.Item(96,3).Value = 8450689526
$ws.Cells.Item(96,4).NumberFormat = "@"
$ws.Cells.Item(96,4).Value = "13052054965"
$ws.Cells.Item(96,4).Style = "Normal"
$ws.Cells.Item(96,5).Value = "Test"

# Row 97
$ws.Cells.Item(97,1).Value = "2025-10-20 19:10:32"
$ws.Cells.Item(97,2).Value = "Noah Dubitzky"
$ws.Cells.Item(97,3).Value = 8450689526
$ws.Cells.Item(97,4).NumberFormat = "@"
$ws.Cells.Item(97,4).Value = "13052054965"
$ws.Cells.Item(97,4).Style = "Normal"
$ws.Cells.Item(97,5).Value = "Test"

# Row 98
$ws.Cells.Item(98,1).Value = "2025-10-20 19:12:10"
$ws.Cells.Item(98,2).Value = "Noah Dubitzky"
$ws.Cells.Item(98,3).Value = 8450689526
$ws.Cells.Item(98,4).NumberFormat = "@"
$ws.Cells.Item(98,4).Value = "13052054965"
$ws.Cells.Item(98,4).Style = "Normal"
$ws.Cells.Item(98,5).Value = "Test"

# Row 99
$ws.Cells.Item(99,1).Value = "2025-10-20 20:25:40"
$ws.Cells.Item(99,2).Value = "Noahs channel"
$ws.Cells.Item(99,4).Value = "Unknown"
$ws.Cells.Item(99,5).Value = "Test"
$ws.Cells.Item(99,7).Value = "Noahs channel"

